# Applies the FlashScore odds updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G3"  = 2.7
    "I3"  = 3
    "M3"  = 1.13
    "N3"  = 6
    "AJ3" = 29
    "AN3" = 4.5

    "G4"  = 2.05
    "I4"  = 4
    "W4"  = 5.5
    "AH4" = 19

    "G14"  = 2.87
    "H14"  = 2.9
    "I14"  = 2.52
    "J14"  = 3.35
    "L14"  = 3.15
    "N14"  = 6.95
    "O14"  = 1.38
    "Q14"  = 2.1
    "R14"  = 1.57
    "W14"  = 8.25
    "X14"  = 14.5
    "Y14"  = 10.25
    "Z14"  = 37
    "AA14" = 25
    "AB14" = 35
    "AC14" = 7.5
    "AG14" = 7
    "AH14" = 11.75
    "AI14" = 9.75
    "AJ14" = 29
    "AK14" = 24
    "AN14" = 4.7
    "AO14" = 15.5
    "AP14" = 22
    "AQ14" = 70
    "AR14" = 100
    "AT14" = 2.4
    "AW14" = 4.35
    "AX14" = 14
    "AY14" = 22
    "AZ14" = 65
    "BA14" = 100
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
